# Generate Report for Handoff
# - Status moves from "Handed back: in sync with en-US" to "Ready for handoff"
# - Per-language "Latest Handoff Datetime" values are refreshed
# - The now-narrower "Status" / language columns are resized to fit the new text

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update status text (same shared text everywhere it appears) ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Refresh the handoff timestamps for each language ---
$zhcn.Range("H2").Value = "2016-09-01 03:10:34"
$dede.Range("H2").Value = "2016-09-01 03:10:38"

# --- Roll the per-language "Latest HO Xliff Generate Date" up to the Overview sheet ---
$overview.Range("G2").Value = "2016-09-01 03:10:38"

# --- Shrink the columns that previously held the long status text ---
$overview.Columns.Item(5).ColumnWidth = 16.3
$overview.Columns.Item(6).ColumnWidth = 16.3
$zhcn.Columns.Item(3).ColumnWidth = 16.3
$dede.Columns.Item(3).ColumnWidth = 16.3
